# Updates the cryptos list values (prices & hourly volume %) per the commit diff.
# Also swaps the Stellar/Filecoin rows (31/32) content as captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    # Force the cell to remain plain text so values such as "0.9400" or
    # "14.80" are not silently re-interpreted (and truncated) as numbers.
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "25.606.89"
Set-TextValue $ws.Range("E2") "  +2.64%  "
Set-TextValue $ws.Range("D3") "1.667.36"
Set-TextValue $ws.Range("E3") "  +1.55%  "
Set-TextValue $ws.Range("D4") "0.9977"
Set-TextValue $ws.Range("E4") "  -0.16%  "
Set-TextValue $ws.Range("D5") "237.43"
Set-TextValue $ws.Range("E5") "  +2.04%  "
Set-TextValue $ws.Range("D6") "0.9995"
Set-TextValue $ws.Range("E6") "  -0.10%  "
Set-TextValue $ws.Range("D7") "0.4622"
Set-TextValue $ws.Range("E7") "  -2.66%  "
Set-TextValue $ws.Range("D8") "0.2576"
Set-TextValue $ws.Range("E8") "  -0.23%  "
Set-TextValue $ws.Range("D9") "0.06135"
Set-TextValue $ws.Range("E9") "  +0.59%  "
Set-TextValue $ws.Range("D10") "1.663.68"
Set-TextValue $ws.Range("E10") "  +1.28%  "
Set-TextValue $ws.Range("D11") "0.06929"
Set-TextValue $ws.Range("E11") "  -1.51%  "
Set-TextValue $ws.Range("D12") "14.87"
Set-TextValue $ws.Range("E12") "  +2.58%  "
Set-TextValue $ws.Range("D13") "4.338"
Set-TextValue $ws.Range("E13") "  +0.07%  "
Set-TextValue $ws.Range("D14") "75.17"
Set-TextValue $ws.Range("E14") "  +2.08%  "
Set-TextValue $ws.Range("D15") "0.5739"
Set-TextValue $ws.Range("E15") "  -2.37%  "
Set-TextValue $ws.Range("E16") "  -0.03%  "
Set-TextValue $ws.Range("D17") "0.9996"
Set-TextValue $ws.Range("D18") "25.571.32"
Set-TextValue $ws.Range("E18") "  +2.48%  "
Set-TextValue $ws.Range("D19") "0.000006685"
Set-TextValue $ws.Range("E19") "  +1.63%  "
Set-TextValue $ws.Range("D20") "11.36"
Set-TextValue $ws.Range("E20") "  +1.39%  "
Set-TextValue $ws.Range("D21") "1.875.36"
Set-TextValue $ws.Range("E21") "  +1.06%  "
Set-TextValue $ws.Range("D22") "4.427"
Set-TextValue $ws.Range("E22") "  +3.25%  "
Set-TextValue $ws.Range("D23") "8.604"
Set-TextValue $ws.Range("E23") "  +0.75%  "
Set-TextValue $ws.Range("D24") "5.227"
Set-TextValue $ws.Range("E24") "  +0.03%  "
Set-TextValue $ws.Range("D25") "134.12"
Set-TextValue $ws.Range("E25") "  +0.23%  "
Set-TextValue $ws.Range("D26") "14.91"
Set-TextValue $ws.Range("E26") "  +0.14%  "
Set-TextValue $ws.Range("E27") "  -0.28%  "
Set-TextValue $ws.Range("D29") "103.85"
Set-TextValue $ws.Range("E29") "  +0.60%  "
Set-TextValue $ws.Range("D30") "3.943"
Set-TextValue $ws.Range("E30") "  +1.77%  "
# Row 31/32: the "Stellar" and "Filecoin" rows swap position (Filecoin moves
# up to row 31, Stellar moves down to row 32), each carrying fresh
# price/volume figures.
Set-TextValue $ws.Range("B31") "Filecoin"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D31") "3.600"
Set-TextValue $ws.Range("E31") "  +0.89%  "

Set-TextValue $ws.Range("B32") "Stellar"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D32") "0.07638"
Set-TextValue $ws.Range("E32") "  +0.65%  "

Set-TextValue $ws.Range("D33") "0.04336"
Set-TextValue $ws.Range("E33") "  +1.41%  "
Set-TextValue $ws.Range("D34") "2.605"
Set-TextValue $ws.Range("E34") "  +1.40%  "
Set-TextValue $ws.Range("D35") "0.6074"
Set-TextValue $ws.Range("E35") "  +3.10%  "
Set-TextValue $ws.Range("D36") "0.9400"
Set-TextValue $ws.Range("E36") "  +1.67%  "
Set-TextValue $ws.Range("D37") "0.9282"
Set-TextValue $ws.Range("E37") "  +6.39%  "
Set-TextValue $ws.Range("D38") "2.462"
Set-TextValue $ws.Range("E38") "  -4.39%  "
Set-TextValue $ws.Range("D39") "106.98"
Set-TextValue $ws.Range("E39") "  +8.42%  "
Set-TextValue $ws.Range("D40") "0.9986"
Set-TextValue $ws.Range("E40") "  -0.15%  "
Set-TextValue $ws.Range("E41") "  +4.10%  "
Set-TextValue $ws.Range("D42") "0.01445"
Set-TextValue $ws.Range("E42") "  -3.75%  "
Set-TextValue $ws.Range("D43") "5.048"
Set-TextValue $ws.Range("E43") "  +8.09%  "
Set-TextValue $ws.Range("D44") "0.3704"
Set-TextValue $ws.Range("E44") "  +0.25%  "
Set-TextValue $ws.Range("D45") "0.1110"
Set-TextValue $ws.Range("E45") "  +0.73%  "
Set-TextValue $ws.Range("D46") "0.05268"
Set-TextValue $ws.Range("E46") "  +1.28%  "
Set-TextValue $ws.Range("D47") "31.21"
Set-TextValue $ws.Range("E47") "  +8.62%  "
Set-TextValue $ws.Range("D48") "6.088"
Set-TextValue $ws.Range("E48") "  +0.04%  "
Set-TextValue $ws.Range("D49") "7.595"
Set-TextValue $ws.Range("E49") "  +7.03%  "
Set-TextValue $ws.Range("E50") "  +0.07%  "
Set-TextValue $ws.Range("D51") "0.9977"
Set-TextValue $ws.Range("E51") "  -0.13%  "
